# tests updated on 30/08/2018
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contact")

# Update the Title values to include a trailing period
$ws.Range("A2").Value = "Mr."
$ws.Range("A3").Value = "Dr."
$ws.Range("A4").Value = "Mrs."

# Move/update the active selection to A4
$ws.Activate()
$ws.Range("A4").Select()
